# BOM updated to Pico2 board with RP2350
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Row 19: swap the RP Pico (with headers) line for the RP Pico2 (RP2350) line.
$ws.Range("F19").Value = "SC1631"
$ws.Range("D19").Value = "2648-SC1631CT-ND"
$ws.Range("B19").Value = "RP Pico2 with RP2350"
$ws.Range("H19").Value = "Pico2 RP2350 CPU with Flash"
$ws.Range("I19").Value = "this version does not have headers"
$ws.Range("C19").Value = "Digikey"
$ws.Range("E19").Value = "Raspberry Pi"
$ws.Range("G19").Value = 1

# Row 20: new line - the Pico2 no longer ships with headers, so add them separately.
$ws.Range("I20").Value = "needed for Pico2 as this version does not have headers"
$ws.Range("H20").Value = "Pico2 headers"
$ws.Range("B20").Value = "Pinstrip male 20 pin"
$ws.Range("E20").Value = "generic"
$ws.Range("G20").Value = 2

# Reflect the new last-used selection on the sheet.
$ws.Range("H20").Select()
